$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.377.27'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '1.649.88'
$ws.Range("E3").Value = '  -2.67%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").Value = "'213.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").Value = "'24.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").Value = "'0.0875"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("D12").Value = '1.880.71'
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("D13").Value = '1.644.37'
$ws.Range("E13").Value = '  -3.10%  '
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("D15").Value = "'0.568"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").Value = "'65.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").Value = '27.363.97'
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").Value = "'234.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.88%  '
$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").Value = "'7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.96%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("D23").Value = "'9.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").Value = "'145.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("D26").Value = "'7.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.16%  '
$ws.Range("D27").Value = "'16.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("E30").Value = '  -2.04%  '
$ws.Range("E31").Value = '  -1.63%  '
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("D33").Value = '1.461.44'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = "'3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("D38").Value = "'0.572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.10%  '
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = "'5.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.11%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = "'65.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("D45").Value = '1.790.10'
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("D46").Value = "'0.780"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").Value = "'88.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("E51").Value = '  -2.50%  '
